# Update the "want to go" counts (column F) for a handful of events on both
# the "展览" (exhibition) sheet and the "全部类型" (all types) sheet, which
# mirror each other's data.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    4  = 1552
    7  = 11254
    10 = 418
    14 = 12288
    15 = 12915
    22 = 73
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
